# Update the "想去人数" (column F) figures on the two sheets that hold
# the full event listing: "展览" and "全部类型". Both sheets mirror the
# same rows, so the same updates are applied to each.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F (想去人数)
$updates = @{
    2  = 1764
    3  = 808
    4  = 4
    7  = 12157
    10 = 485
    11 = 423
    12 = 1122
    13 = 880
    14 = 13551
    15 = 13635
    17 = 158
    20 = 1008
    23 = 2079
    24 = 192
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
